$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (179) down to the new rows (180:191)
$ws.Range("A179:I179").Copy()
$ws.Range("A180:I191").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# For rows where column G (Localisation douleur) holds actual text, the formatting
# differs from the blank-G style (centered Helvetica) used by row 179: it uses the
# same left-aligned "Helvetica Neue" style as the other text columns (e.g. column B).
$ws.Range("B179").Copy()
$ws.Range("G184").PasteSpecial(-4122)
$ws.Range("G187").PasteSpecial(-4122)
$ws.Range("G188").PasteSpecial(-4122)
$ws.Range("G189").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 180
$ws.Range("A180").Value = 45890
$ws.Range("B180").Value = "Ilyes Boughanmi"
$ws.Range("C180").Value = 75
$ws.Range("D180").Value = 6
$ws.Range("E180").Value = 6
$ws.Range("F180").Value = 0
$ws.Range("H180").Value = 0
$ws.Range("I180").Formula = "=C180*D180"

# Row 181
$ws.Range("A181").Value = 45890
$ws.Range("B181").Value = "Omar Benyounes"
$ws.Range("C181").Value = 75
$ws.Range("D181").Value = 6
$ws.Range("E181").Value = 5
$ws.Range("F181").Value = 0
$ws.Range("H181").Value = 5
$ws.Range("I181").Formula = "=C181*D181"

# Row 182
$ws.Range("A182").Value = 45890
$ws.Range("B182").Value = "Naim Ighbane"
$ws.Range("C182").Value = 75
$ws.Range("D182").Value = 6
$ws.Range("E182").Value = 0
$ws.Range("F182").Value = 0
$ws.Range("H182").Value = 0
$ws.Range("I182").Formula = "=C182*D182"

# Row 183
$ws.Range("A183").Value = 45890
$ws.Range("B183").Value = "Yanis Berrached"
$ws.Range("C183").Value = 75
$ws.Range("D183").Value = 7
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 0
$ws.Range("H183").Value = 6
$ws.Range("I183").Formula = "=C183*D183"

# Row 184
$ws.Range("A184").Value = 45890
$ws.Range("B184").Value = "Maé Clavel"
$ws.Range("C184").Value = 75
$ws.Range("D184").Value = 7
$ws.Range("E184").Value = 5
$ws.Range("F184").Value = 3
$ws.Range("G184").Value = "Cheville"
$ws.Range("H184").Value = 7
$ws.Range("I184").Formula = "=C184*D184"

# Row 185
$ws.Range("A185").Value = 45890
$ws.Range("B185").Value = "Romain Thunet"
$ws.Range("C185").Value = 75
$ws.Range("D185").Value = 5
$ws.Range("E185").Value = 3
$ws.Range("F185").Value = 0
$ws.Range("H185").Value = 2
$ws.Range("I185").Formula = "=C185*D185"

# Row 186
$ws.Range("A186").Value = 45890
$ws.Range("B186").Value = "Ilan Ihaddadene"
$ws.Range("C186").Value = 75
$ws.Range("D186").Value = 7
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 0
$ws.Range("H186").Value = 9
$ws.Range("I186").Formula = "=C186*D186"

# Row 187
$ws.Range("A187").Value = 45890
$ws.Range("B187").Value = "Naim Dhib"
$ws.Range("C187").Value = 75
$ws.Range("D187").Value = 5
$ws.Range("E187").Value = 4
$ws.Range("F187").Value = 3
$ws.Range("G187").Value = "Aine"
$ws.Range("H187").Value = 6
$ws.Range("I187").Formula = "=C187*D187"

# Row 188
$ws.Range("A188").Value = 45890
$ws.Range("B188").Value = "Levy Ndoutoume"
$ws.Range("C188").Value = 75
$ws.Range("D188").Value = 7
$ws.Range("E188").Value = 7
$ws.Range("F188").Value = 5
$ws.Range("G188").Value = "Ischio"
$ws.Range("H188").Value = 7
$ws.Range("I188").Formula = "=C188*D188"

# Row 189
$ws.Range("A189").Value = 45890
$ws.Range("B189").Value = "Karahali Souaré"
$ws.Range("C189").Value = 75
$ws.Range("D189").Value = 5
$ws.Range("E189").Value = 6
$ws.Range("F189").Value = 7
$ws.Range("G189").Value = "Creux poplité"
$ws.Range("H189").Value = 5
$ws.Range("I189").Formula = "=C189*D189"

# Row 190
$ws.Range("A190").Value = 45890
$ws.Range("B190").Value = "Mattheo Haon"
$ws.Range("C190").Value = 75
$ws.Range("D190").Value = 5
$ws.Range("E190").Value = 0
$ws.Range("F190").Value = 0
$ws.Range("H190").Value = 4
$ws.Range("I190").Formula = "=C190*D190"

# Row 191
$ws.Range("A191").Value = 45890
$ws.Range("B191").Value = "Sofiane Belle"
$ws.Range("C191").Value = 75
$ws.Range("D191").Value = 5
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 0
$ws.Range("H191").Value = 6
$ws.Range("I191").Formula = "=C191*D191"

# Update sheet view to match the final state (selection)
$ws.Range("I184").Select()